$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# Row 15 corresponds to "crude oil" electricity source.
# Set its BAU Guaranteed Dispatch Percentage to 0 for all year columns (B:AK, 2015-2050),
# reflecting that crude oil / cogeneration dispatch is no longer guaranteed.
$ws.Range("B15:AK15").Value = 0
